# Adds the 14 Mayis 2020 (2020-05-14) row of data to the "data" worksheet
# and lets the existing Excel Table (Table3) grow to include it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Copy formatting from the last existing data row (63) down to the new row (64)
# so the new cells reuse the same cell styles (date format on column A, etc.)
$ws.Range("A63:E63").Copy()
$ws.Range("A64:E64").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New row of data (row 64): date, test, case, death, recovered
$ws.Range("A64").Value = 43965
$ws.Range("B64").Value = 34821
$ws.Range("C64").Value = 1635
$ws.Range("D64").Value = 55
$ws.Range("E64").Value = 2315

# Grow the table (ListObject) so its range / autofilter include the new row
$table = $ws.ListObjects.Item("Table3")
$table.Resize($ws.Range("A1:E64"))

$wb.Save()
